# Update TPM-derived metrics in the LR-pairs sheet (Ybx1-Notch1)
# Columns G:T for rows 2-10 are recalculated values coming from the
# updated TPM input data (see commit message "update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value
$updates = @{
    2 = @{ G = 175.9726616666667;  H = 527.917985;           I = 0.5957262918830742;  J = 0.5957262918830741;
           M = 48.42420966666666;  N = 145.272629;           O = 0.6311762527593259;  P = 0.6311762527593258;
           Q = 8521.337064148063;  R = 76692.03357733256;    S = 0.3760082885809672;  T = 0.3760082885809671 }
    3 = @{ G = 175.9726616666667;  H = 527.917985;           I = 0.5957262918830742;  J = 0.5957262918830741;
           M = 6.849914666666667;  N = 20.549744;            O = 0.08928392431779728; P = 0.08928392431779726;
           Q = 1205.397716082871;  R = 10848.57944474584;    S = 0.05318878115861041; T = 0.05318878115861039 }
    4 = @{ G = 175.9726616666667;  H = 527.917985;           I = 0.5957262918830742;  J = 0.5957262918830741;
           M = 21.446458;          N = 64.33937399999999;    O = 0.2795398229228769;  P = 0.2795398229228769;
           Q = 3773.990297582376;  R = 33965.91267824139;    S = 0.1665292221434966;  T = 0.1665292221434966 }
    5 = @{ I = 0.2361302886539166; J = 0.2361302886539166;
           M = 48.42420966666666;  N = 145.272629;           O = 0.6311762527593259;  P = 0.6311762527593258;
           Q = 3377.634675673393;  R = 30398.71208106053;    S = 0.149039830755557;   T = 0.149039830755557 }
    6 = @{ I = 0.2361302886539166; J = 0.2361302886539166;
           M = 6.849914666666667;  O = 0.08928392431779728;  P = 0.08928392431779726;
           S = 0.02108263882131592; T = 0.02108263882131591 }
    7 = @{ I = 0.2361302886539166; J = 0.2361302886539166;
           N = 64.33937399999999;  O = 0.2795398229228769;   P = 0.2795398229228769;
           S = 0.06600781907704366; T = 0.06600781907704363 }
    8 = @{ G = 49.66818733333334;  I = 0.1681434194630093;   J = 0.1681434194630093;
           M = 48.42420966666666;  N = 145.272629;           O = 0.6311762527593259;  P = 0.6311762527593258;
           Q = 2405.142717192611;  R = 21646.2844547335;     S = 0.1061281334228017;  T = 0.1061281334228017 }
    9 = @{ G = 49.66818733333334;  I = 0.1681434194630093;   J = 0.1681434194630093;
           M = 6.849914666666667;  O = 0.08928392431779728;  P = 0.08928392431779726;
           Q = 340.2228448813476;  S = 0.01501250433787096;  T = 0.01501250433787096 }
    10 = @{ G = 49.66818733333334; I = 0.1681434194630093;   J = 0.1681434194630093;
           N = 64.33937399999999;  O = 0.2795398229228769;   P = 0.2795398229228769;
           R = 9586.860242224187;  S = 0.04700278170233663;  T = 0.04700278170233662 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
